# Update Smith's phone number (C8) to a full international number and
# make sure it is stored/formatted as text so the leading "+" survives.
$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws1.Range("C8").NumberFormat = "@"
$ws1.Range("C8").Value = "+256790513253"

# Rows 8-9 picked up the tighter default row height that comes with the
# re-saved workbook.
$ws1.Rows.Item(8).RowHeight = 13.8
$ws1.Rows.Item(9).RowHeight = 13.8

# Leave the UI focused on Sheet1 with C11:C12 selected, as it was when the
# edit was made.
$ws1.Activate()
$null = $ws1.Range("C11:C12").Select()
